$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 63.211268
$ws.Range("H2").Value = 189.633804
$ws.Range("I2").Value = 0.4922609885657722
$ws.Range("J2").Value = 0.4922609885657722
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8412133333333333
$ws.Range("N2").Value = 2.52364
$ws.Range("Q2").Value = 53.17416145850666
$ws.Range("R2").Value = 478.56745312656
$ws.Range("S2").Value = 0.4922609885657722
$ws.Range("T2").Value = 0.4922609885657722

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 43.30706799999999
$ws.Range("H3").Value = 129.921204
$ws.Range("I3").Value = 0.3372560111523963
$ws.Range("J3").Value = 0.3372560111523963
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.8412133333333333
$ws.Range("N3").Value = 2.52364
$ws.Range("Q3").Value = 36.43048302917332
$ws.Range("R3").Value = 327.87434726256
$ws.Range("S3").Value = 0.3372560111523963
$ws.Range("T3").Value = 0.3372560111523963

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.891734
$ws.Range("H4").Value = 65.675202
$ws.Range("I4").Value = 0.1704830002818315
$ws.Range("J4").Value = 0.1704830002818315
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.8412133333333333
$ws.Range("N4").Value = 2.52364
$ws.Range("Q4").Value = 18.41561853058667
$ws.Range("R4").Value = 165.74056677528
$ws.Range("S4").Value = 0.1704830002818315
$ws.Range("T4").Value = 0.1704830002818315
